$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.198.86"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.930.55"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'592.97"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "'145.76"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.507"
$ws.Range("E8").Value = "  +1.21%  "
$ws.Range("D9").Value = "2.929.87"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").Value = "'0.145"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'0.444"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "'33.79"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "3.412.72"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "61.080.42"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "'6.75"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "2.918.97"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'432.94"
$ws.Range("E20").Value = "  +1.14%  "
$ws.Range("D21").Value = "'13.52"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").Value = "'0.684"
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("D23").Value = "'7.11"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "'80.99"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'10.86"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").Value = "'2.23"
$ws.Range("E26").Value = "  +2.04%  "
$ws.Range("D27").Value = "'12.19"
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").Value = "'2.35"
$ws.Range("E29").Value = "  +7.20%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'7.13"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "'26.67"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("E35").Value = "  +4.92%  "
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.12"
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'5.65"
$ws.Range("E38").Value = "  -0.96%  "
$ws.Range("D39").Value = "'50.10"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'2.02"
$ws.Range("E40").Value = "  -0.27%  "
$ws.Range("E41").Value = "  -1.58%  "
$ws.Range("D42").Value = "'8.64"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Value = "'0.290"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'39.98"
$ws.Range("E44").Value = "  -4.17%  "
$ws.Range("D45").Value = "'384.20"
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("D46").Value = "'0.0350"
$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("D47").Value = "2.709.87"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").Value = "'129.84"
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.38"
$ws.Range("E50").Value = "  -3.89%  "
$ws.Range("E51").Value = "  +0.57%  "
